$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data was collected for Cilantro at Terminal La Palmera
# de La Serena. Insert a new row above the current first data row (row 33)
# so the existing rows 33-68 shift down to 34-69, then populate the new
# row 33 with this week's figures.
$ws.Rows(33).Insert()

$ws.Cells.Item(33, 1).Value = 8
$ws.Cells.Item(33, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(33, 3).Value = "Coquimbo"
$ws.Cells.Item(33, 4).Value = 44413
$ws.Cells.Item(33, 5).Value = 4
$ws.Cells.Item(33, 6).Value = 100112040
$ws.Cells.Item(33, 7).Value = "Cilantro"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 3600
$ws.Cells.Item(33, 11).Value = 2000
$ws.Cells.Item(33, 12).Value = 2500
$ws.Cells.Item(33, 13).Value = 2250
$ws.Cells.Item(33, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(33, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(33, 16).Value = 1500
$ws.Cells.Item(33, 17).Value = 1.5
$ws.Cells.Item(33, 18).Value = "Hortaliza"
